$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.852.19'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '3.839.64'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.25'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.35'
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('D7').Value = '3.833.26'
$ws.Range('E7').Value = '  +2.04%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.41'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000262'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.68'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '4.479.94'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '3.846.40'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '68.868.61'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.18'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.05'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.112'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.98'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '469.03'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000156'
$ws.Range('E24').Value = '  +10.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.56'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').Value = '  -2.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.01'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.22'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.79'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.41'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.40'
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.20'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.24'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.793.43'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.102'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.50'
$ws.Range('E38').Value = '  +2.36%  '
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.84'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.20'
$ws.Range('E44').Value = '  +16.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.305'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.95'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.90'
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.49'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '146.89'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '394.54'
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('D51').Value = '2.812.03'
$ws.Range('E51').Value = '  +4.64%  '
